# Update predicted values in column C (rows 2-201) of the active worksheet
# with the new model predictions from the latest run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    72.44926452636719,180.0577545166016,121.7747116088867,37.84516143798828,16.25225639343262,13.2372932434082,17.38755989074707,32.64756393432617,90.34137725830078,139.6136932373047,
    66.93205261230469,22.57066345214844,13.33227348327637,15.24962615966797,28.31535339355469,61.92763137817383,113.8759841918945,77.08506774902344,36.45169067382812,23.05905532836914,
    19.89080619812012,33.86481857299805,69.91464996337891,116.5624313354492,91.73889923095703,42.60548400878906,25.61328887939453,26.48396110534668,38.65092468261719,70.02791595458984,
    108.7448348999023,92.27141571044922,46.70959091186523,26.11845397949219,24.41765403747559,35.98159790039062,69.67458343505859,115.8239517211914,98.66626739501953,49.67593002319336,
    30.3542423248291,29.6086483001709,38.18948745727539,66.18819427490234,108.7140045166016,99.96614837646484,55.13914108276367,33.41853713989258,29.98615455627441,39.29381942749023,
    66.83805847167969,110.0486297607422,102.4283218383789,57.81891250610352,38.3359375,35.29103469848633,42.52236557006836,68.65338897705078,107.1457901000977,102.9799880981445,
    62.48945999145508,41.88785934448242,37.65022659301758,46.14177322387695,72.26847839355469,107.2397232055664,102.8035430908203,64.79190063476562,46.33440780639648,42.55091857910156,
    50.49316787719727,76.55509948730469,106.4182891845703,102.5270385742188,68.06403350830078,50.08418655395508,46.69758224487305,56.17034530639648,82.78958892822266,108.4250564575195,
    102.369743347168,70.50347900390625,54.24306106567383,51.70468521118164,62.09638595581055,89.92745208740234,111.9324417114258,103.2621536254883,73.1793212890625,57.57878494262695,
    55.35265731811523,67.19908905029297,96.35614776611328,115.7222595214844,104.2056884765625,75.24977874755859,60.93056106567383,58.13800811767578,70.29374694824219,100.3155364990234,
    119.3992691040039,106.4768371582031,77.36252593994141,63.85659408569336,61.94568252563477,75.13082122802734,105.7200393676758,122.9576416015625,108.6207962036133,80.53170013427734,
    67.86582946777344,67.13919830322266,82.39321136474609,115.1204986572266,129.3532867431641,111.557258605957,83.97867584228516,71.77127838134766,72.55974578857422,90.23355102539062,
    125.9117813110352,136.8389892578125,114.4787673950195,87.14205169677734,74.98500823974609,77.04271697998047,97.17916107177734,133.6876373291016,143.0653228759766,117.6526565551758,
    90.29673767089844,78.7705078125,82.54402923583984,105.0641403198242,141.7783660888672,149.2194366455078,121.2257308959961,94.52275085449219,84.12854766845703,89.95444488525391,
    115.6325302124023,150.5804443359375,154.8553924560547,125.3143615722656,100.0630722045898,91.78626251220703,100.9035415649414,130.7525634765625,161.5493621826172,160.5995025634766,
    130.0408630371094,108.0814361572266,102.2063140869141,115.554801940918,151.0543518066406,175.9818725585938,167.2915496826172,135.2693023681641,118.4857482910156,118.5530471801758,
    139.2193298339844,177.0422821044922,191.4904022216797,174.4512481689453,144.4447174072266,134.8780822753906,144.0814666748047,176.6115264892578,203.5403900146484,203.9424743652344,
    183.2945709228516,162.1428680419922,161.4267730712891,183.8692779541016,213.2490844726562,227.3411407470703,219.7437896728516,200.5043487548828,189.544677734375,200.7989654541016,
    225.9990234375,247.4006042480469,252.0785980224609,240.5500640869141,226.4314422607422,226.1552124023438,242.9902801513672,266.7535095214844,281.7860717773438,279.6876831054688,
    267.40087890625,260.50048828125,269.1564025878906,290.1633911132812,310.0154724121094,316.8477783203125,310.9211120605469,303.1871337890625,304.9046325683594,319.2162780761719
)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $rowIndex = $i + 2
    $ws.Cells.Item($rowIndex, 3).Value = [double]$newValues[$i]
}
